$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the data ranges as Text so that numeric-looking strings
# (e.g. "1.000", "0.07883") are not silently coerced into numbers by Excel.
$dRange = $ws.Range("D2:D51")
$eRange = $ws.Range("E2:E51")
$dRange.NumberFormat = "@"
$eRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.192.01"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.852.07"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").Value = "237.56"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.07883"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").Value = "0.3013"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "23.67"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").Value = "0.08102"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.855.10"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "5.174"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "0.7038"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("D15").Value = "89.44"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "29.226.44"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "5.798"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "0.000007812"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "235.45"
$ws.Range("D22").Value = "2.101.52"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "7.489"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "162.44"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").Value = "8.854"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "18.00"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "1.404"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "1.476"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("D33").Value = "4.009"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "0.05154"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "1.162"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").Value = "0.7100"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "0.9985"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").Value = "2.679"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "0.01844"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "2.705"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("D41").Value = "1.149.65"
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("D42").Value = "0.9219"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "0.4229"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").Value = "69.95"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "102.89"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "0.5292"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("D49").Value = "1.734"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").Value = "9.137"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  -0.74%  "

# Restore the default (Normal) style so no stray formatting is left behind.
$dRange.Style = "Normal"
$eRange.Style = "Normal"
